$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab to the new ASIN code.
$ws.Name = "B08K4353V1"

# Replace the 100 keyword rows in column A with the new keyword list, in order.
$values = @(
    'neck compression wrap',
    'training bra sets',
    'wide waistband yoga',
    'girls seamless underwear',
    'sport bra',
    'seamless bra',
    'yoga legging',
    'workout waist',
    'yoga workout',
    'yoga bra',
    'sports bra for women',
    'sport apparel',
    'seamless bra set',
    'workout apparel',
    'workout bra',
    'workout outfit set',
    'workout set',
    'workout leggings for women high waist',
    'workout legging',
    'legging bra set',
    'workout outfit',
    'sports leggings women',
    'black legging',
    'sport bra set',
    'yoga sets for women',
    'high waist workout',
    'yoga set purple',
    'yoga bras for women',
    '2 pcs',
    'green 2',
    'pcs set',
    'bra set',
    'waist yoga',
    'outfit set',
    'seamless sports bra',
    'seamless workout leggings',
    'seamless yoga bra',
    'seamless yoga',
    'sports for women',
    'high waist yoga',
    'sport legging',
    'black bra',
    'yoga legging set',
    'yoga workout set',
    'yoga bra set',
    'yoga outfit set',
    'workout set women',
    'black legging set',
    'sports leggings for women',
    'black bra set',
    'sport leggings for women',
    'sport set',
    'black outfit',
    'high waist',
    'set 2',
    'yoga apparel',
    'purple bra',
    'purple legging',
    'high waist legging',
    'high sport bra',
    'black n',
    'womens sports bra set',
    'sport outfit',
    'sport waist',
    'sport workout',
    'waist bra',
    'black 2',
    'yoga sport',
    'n set',
    'green set',
    'yoga 2',
    'green bra',
    'green bra set',
    'green yoga',
    'yoga sport bra',
    'womens workout sports bra',
    'high waist set',
    'ready set',
    'high black',
    'waist set',
    'black yoga',
    'green outfit',
    'black apparel',
    'purple yoga',
    'black set',
    'purple outfit',
    'purple workout',
    'black sport',
    'seamless legging',
    'green apparel',
    'purple set',
    'black seamless bra',
    'black yoga bra',
    'seamless workout',
    'green sport bra',
    'black sport bra',
    'purple sport bra',
    'workout sport bra',
    'black workout bra',
    'seamless sport bra'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

# Update the sheet view selection to the full data range.
$ws.Range("A1:A100").Select() | Out-Null
